# "Generate Report for Archive"
#
# The status text "Ready for handoff" has changed to "In Translation"
# everywhere it is used (the Overview sheet's per-locale status columns,
# plus each locale sheet's own Status column). Because the new text is
# shorter than the old text, the (previously auto-fitted) status columns
# need to shrink to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status value wherever it appears.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value     = "In Translation"
$dede.Range("C2").Value     = "In Translation"

# Re-fit the columns that held the status text so their width tracks the
# new (shorter) content, same as Excel would do for an auto-fitted column.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth     = 12.5
$dede.Columns.Item(3).ColumnWidth     = 12.5
